$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values for columns D (Japon - todos los items) and E (Japon - menos
# alimentos y energia) for rows 287 through 307 (revised figures), plus a
# brand-new monthly record appended as row 308.
$updates = @(
    @{ Row = 287; D = 100.4;              E = 100.4 },
    @{ Row = 288; D = 100.5;              E = 100.4 },
    @{ Row = 289; D = 100.5;              E = 100.4 },
    @{ Row = 290; D = 100.5;              E = 100.2 },
    @{ Row = 291; D = 100.3;              E = 100.1 },
    @{ Row = 292; D = 100.3;              E = 100.3 },
    @{ Row = 293; D = 100.2;              E = 100.1 },
    @{ Row = 294; D = 100.1;              E = 100.2 },
    @{ Row = 295; D = 99.90000000000001;  E = 100 },
    @{ Row = 296; D = 100;                E = 100 },
    @{ Row = 297; D = 100.1;              E = 99.8 },
    @{ Row = 298; D = 99.90000000000001;  E = 99.8 },
    @{ Row = 299; D = 99.8;               E = 99.90000000000001 },
    @{ Row = 300; D = 99.5;               E = 99.90000000000001 },
    @{ Row = 301; D = 99.3;               E = 99.90000000000001 },
    @{ Row = 302; D = 99.8;               E = 100.2 },
    @{ Row = 303; D = 99.8;               E = 100.2 },
    @{ Row = 304; D = 99.90000000000001;  E = 100.3 },
    @{ Row = 305; D = 99.09999999999999;  E = 98.90000000000001 },
    @{ Row = 306; D = 99.40000000000001;  E = 99 },
    @{ Row = 307; D = 99.5;               E = 98.90000000000001 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Append the new monthly record as row 308.
# The date label "01-07-2021" is built through a scratch formula cell and
# copy/pasted in as a value so the COM layer stores it as plain text (matching
# the existing text-based date labels) instead of auto-converting it to a
# serial date number.
$ws.Range("Z1").Formula = "=""01-07-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A308").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Cells.Item(308, 2).Value = 107.6
$ws.Cells.Item(308, 3).Value = 106.1
$ws.Cells.Item(308, 4).Value = 99.7
$ws.Cells.Item(308, 5).Value = 99.09999999999999
